# Insert a new weekly record for "Vega Modelo de Temuco - Perejil" at row 148.
# This pushes the existing rows 148-210 down to 149-211 (Excel's native
# Insert behavior), and then the newly inserted row 148 is populated with
# the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 148, shifting rows 148:210
# down to 149:211.
$ws.Rows.Item(148).Insert()

# Fill in the new row 148 with the new weekly record. The non-varying
# columns (A,B,C,E,F,G,H,I,N,Q,R) reuse the same constant values found
# throughout this block of rows.
$ws.Cells.Item(148, 1).Value = 10
$ws.Cells.Item(148, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(148, 3).Value = "La Araucanía"
$ws.Cells.Item(148, 4).Value = 44489
$ws.Cells.Item(148, 5).Value = 9
$ws.Cells.Item(148, 6).Value = 100112044
$ws.Cells.Item(148, 7).Value = "Perejil"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 65
$ws.Cells.Item(148, 11).Value = 3000
$ws.Cells.Item(148, 12).Value = 3000
$ws.Cells.Item(148, 13).Value = 3000
$ws.Cells.Item(148, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(148, 15).Value = "Región Metropolitana"
$ws.Cells.Item(148, 16).Value = 1000
$ws.Cells.Item(148, 17).Value = 3
$ws.Cells.Item(148, 18).Value = "Hortaliza"

# Note: Rows.Item(148).Insert() already carries the date number-format
# style (s="2") for the new row's D cell, matching the rest of column D
# in this block, so no extra style assignment is needed here (doing so
# would reset the cell back to the default "Normal" style).
